$d = $word.ActiveDocument

# Color constant for the blue comment-style text (RGB 0x3333FF encoded as BGR integer for Word)
$blue = 16724787

# ---------------------------------------------------------------------------
# 1. "Approval should not be saved..." paragraph: add trailing space + (done)
# ---------------------------------------------------------------------------
$rng = $d.Content
$searchText = "Approval should not be saved if no recommendations or if approved amount is more than the recommended amount, it shall require override and with short explanation why such application was approved beyond the recommended amount."
$found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" ")
    $rng.Collapse(0)
    $startPos = $rng.Start
    $rng.InsertAfter("(done)")
    $newRange = $d.Range($startPos, $startPos + 6)
    $newRange.Font.Color = $blue
}

# ---------------------------------------------------------------------------
# 2. Table: tblInd -45 -> -55, tblCellMar/left 63 -> 53, all tcMar/left 63 -> 53
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.LeftPadding = 2.65
$t.Rows.LeftIndent = -2.75
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.LeftPadding = 2.65
    }
}

# ---------------------------------------------------------------------------
# 3. "(under development)" -> "(" + "done" + ")" (three runs) after the
#    "Maximum loanable amount for two regular account..." paragraph
# ---------------------------------------------------------------------------
$rng = $d.Content
$searchText = "not to exceed the maximum loanable amount for regular account. "
$found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $startPos = $rng.Start
    $targetRng = $d.Range($startPos, $startPos + 20)
    # sanity: $targetRng.Text should be "(under development)"
    $targetRng.Text = ""

    $p1 = $d.Range($startPos, $startPos)
    $p1.InsertAfter("(")
    $p1End = $startPos + 1
    $p1Rng = $d.Range($startPos, $p1End)
    $p1Rng.Font.Color = $blue

    $p2 = $d.Range($p1End, $p1End)
    $p2.InsertAfter("done")
    $p2End = $p1End + 4
    $p2Rng = $d.Range($p1End, $p2End)
    $p2Rng.Font.Color = $blue

    $p3 = $d.Range($p2End, $p2End)
    $p3.InsertAfter(")")
    $p3End = $p2End + 1
    $p3Rng = $d.Range($p2End, $p3End)
    $p3Rng.Font.Color = $blue
}

# ---------------------------------------------------------------------------
# 4. Merge "(" + "is this for reporting purposes? If so, no need to do this." + ")"
#    into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$searchText = "(is this for reporting purposes? If so, no need to do this.)"
$rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $searchText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "Filter lists of competitors by branch in the loan assessment " -> add (done)
# ---------------------------------------------------------------------------
$rng = $d.Content
$searchText = "Filter lists of competitors by branch in the loan assessment "
$found = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $startPos = $rng.Start
    $rng.InsertAfter("(done)")
    $newRange = $d.Range($startPos, $startPos + 6)
    $newRange.Font.Color = $blue
}

# ---------------------------------------------------------------------------
# 6. Bookmark reorder: insert a new bookmark "__DdeLink__275_1647773266"
#    before the existing "__DdeLink__251_845001198" / "__DdeLink__263_46480473"
#    bookmarks (all wrapping the "(done)" after "Show client Net Take Home Pay...")
# ---------------------------------------------------------------------------
$bmOuter = $d.Bookmarks.Item("__DdeLink__263_46480473")
$rngDone = $bmOuter.Range.Duplicate
$bmOuter.Delete()
$bmInner = $d.Bookmarks.Item("__DdeLink__251_845001198")
$bmInner.Delete()

$d.Bookmarks.Add("__DdeLink__275_1647773266", $rngDone)
$d.Bookmarks.Add("__DdeLink__251_845001198", $rngDone)
$d.Bookmarks.Add("__DdeLink__263_46480473", $rngDone)

# ---------------------------------------------------------------------------
# 7. Merge "(" + "what documents?" + ")" into a single run (bookmark id will
#    auto-renumber because of the bookmark inserted in step 6 above).
# ---------------------------------------------------------------------------
$rng = $d.Content
$searchText = "(what documents?)"
$rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $searchText, 2) | Out-Null
